$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2238.5625
$ws.Cells.Item(17, 10).Value = 2238.5625
$ws.Cells.Item(17, 12).Value = 6715.6875
$ws.Cells.Item(17, 14).Value = -7051.6875
$ws.Cells.Item(33, 8).Value = 662.13336
$ws.Cells.Item(33, 10).Value = 1133
$ws.Cells.Item(33, 12).Value = 1133
$ws.Cells.Item(33, 14).Value = -1591
$ws.Cells.Item(40, 8).Value = 11131.211
$ws.Cells.Item(40, 10).Value = 11131.211
$ws.Cells.Item(40, 12).Value = 11131.211
$ws.Cells.Item(40, 14).Value = -11481.211
$ws.Cells.Item(69, 8).Value = 30625
$ws.Cells.Item(69, 9).Value = 30833.334
$ws.Cells.Item(69, 10).Value = 30000
$ws.Cells.Item(69, 11).Value = 92500.00199999999
$ws.Cells.Item(69, 12).Value = 90000
$ws.Cells.Item(69, 13).Value = -91626.00199999999
$ws.Cells.Item(69, 14).Value = -91748
$ws.Cells.Item(72, 8).Value = 30625
$ws.Cells.Item(72, 9).Value = 30833.334
$ws.Cells.Item(72, 10).Value = 30000
$ws.Cells.Item(72, 11).Value = 277500.006
$ws.Cells.Item(72, 12).Value = 270000
$ws.Cells.Item(72, 13).Value = -273132.006
$ws.Cells.Item(72, 14).Value = -278736
$ws.Cells.Item(74, 8).Value = 4006.4614
$ws.Cells.Item(74, 9).Value = 3517.9
$ws.Cells.Item(74, 11).Value = 3517.9
$ws.Cells.Item(74, 13).Value = -2581.9
$ws.Cells.Item(77, 8).Value = 4006.4614
$ws.Cells.Item(77, 9).Value = 3517.9
$ws.Cells.Item(77, 11).Value = 17589.5
$ws.Cells.Item(77, 13).Value = -12909.5
$ws.Cells.Item(109, 8).Value = 91775
$ws.Cells.Item(109, 10).Value = 91775
$ws.Cells.Item(109, 12).Value = 91775
$ws.Cells.Item(109, 14).Value = -94549
$ws.Cells.Item(116, 8).Value = 1048744.8
$ws.Cells.Item(116, 9).Value = 6948.5454
$ws.Cells.Item(116, 11).Value = 6948.5454
$ws.Cells.Item(116, 13).Value = -3506.5454
$ws.Cells.Item(117, 8).Value = 89521.73
$ws.Cells.Item(117, 10).Value = 89521.73
$ws.Cells.Item(117, 12).Value = 89521.73
$ws.Cells.Item(117, 14).Value = -98699.73
$ws.Cells.Item(128, 8).Value = 75000
$ws.Cells.Item(128, 10).Value = 75000
$ws.Cells.Item(128, 12).Value = 75000
$ws.Cells.Item(128, 14).Value = -84960

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8119.4
$ws.Cells.Item(32, 9).Value = 3275.9412
$ws.Cells.Item(32, 11).Value = 3275.9412
$ws.Cells.Item(32, 13).Value = -2988.9412
$ws.Cells.Item(45, 8).Value = 3213.8333
$ws.Cells.Item(45, 10).Value = 3243.75
$ws.Cells.Item(45, 12).Value = 3243.75
$ws.Cells.Item(45, 14).Value = -3997.75
$ws.Cells.Item(61, 8).Value = 2999.8
$ws.Cells.Item(61, 9).Value = 1949.5
$ws.Cells.Item(61, 11).Value = 1949.5
$ws.Cells.Item(61, 13).Value = -1737.5
$ws.Cells.Item(74, 8).Value = 2138.4614
$ws.Cells.Item(74, 9).Value = 1499.75
$ws.Cells.Item(74, 11).Value = 1499.75
$ws.Cells.Item(74, 13).Value = -625.75
$ws.Cells.Item(77, 8).Value = 2138.4614
$ws.Cells.Item(77, 9).Value = 1499.75
$ws.Cells.Item(77, 11).Value = 7498.75
$ws.Cells.Item(77, 13).Value = -3130.75
$ws.Cells.Item(136, 8).Value = 2999.8
$ws.Cells.Item(136, 9).Value = 1949.5
$ws.Cells.Item(136, 11).Value = 5848.5
$ws.Cells.Item(136, 13).Value = -3298.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 6464.6665
$ws.Cells.Item(86, 9).Value = 5697.5
$ws.Cells.Item(86, 10).Value = 7999
$ws.Cells.Item(86, 11).Value = 5697.5
$ws.Cells.Item(86, 12).Value = 7999
$ws.Cells.Item(86, 13).Value = -4574.5
$ws.Cells.Item(86, 14).Value = -10245
$ws.Cells.Item(89, 8).Value = 6464.6665
$ws.Cells.Item(89, 9).Value = 5697.5
$ws.Cells.Item(89, 10).Value = 7999
$ws.Cells.Item(89, 11).Value = 28487.5
$ws.Cells.Item(89, 12).Value = 39995
$ws.Cells.Item(89, 13).Value = -22871.5
$ws.Cells.Item(89, 14).Value = -51227
$ws.Cells.Item(132, 8).Value = 27188.781
$ws.Cells.Item(132, 10).Value = 27188.781
$ws.Cells.Item(132, 12).Value = 27188.781
$ws.Cells.Item(132, 14).Value = -37308.781

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3001.6667
$ws.Cells.Item(16, 9).Value = 2670
$ws.Cells.Item(16, 11).Value = 2670
$ws.Cells.Item(16, 13).Value = -2383
$ws.Cells.Item(31, 8).Value = 12826.194
$ws.Cells.Item(31, 9).Value = 2712.1667
$ws.Cells.Item(31, 11).Value = 2712.1667
$ws.Cells.Item(31, 13).Value = -2417.1667
$ws.Cells.Item(34, 8).Value = 12826.194
$ws.Cells.Item(34, 9).Value = 2712.1667
$ws.Cells.Item(34, 11).Value = 2712.1667
$ws.Cells.Item(34, 13).Value = -2510.1667
$ws.Cells.Item(76, 8).Value = 2385666.2
$ws.Cells.Item(76, 9).Value = 2385666.2
$ws.Cells.Item(76, 11).Value = 2385666.2
$ws.Cells.Item(76, 13).Value = -2385351.2
$ws.Cells.Item(79, 8).Value = 2385666.2
$ws.Cells.Item(79, 9).Value = 2385666.2
$ws.Cells.Item(79, 11).Value = 2385666.2
$ws.Cells.Item(79, 13).Value = -2384574.2
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 3001.6667
$ws.Cells.Item(113, 9).Value = 2670
$ws.Cells.Item(113, 11).Value = 2670
$ws.Cells.Item(113, 13).Value = -500
$ws.Cells.Item(134, 8).Value = 2751.5
$ws.Cells.Item(134, 9).Value = 2278.8
$ws.Cells.Item(134, 11).Value = 6836.400000000001
$ws.Cells.Item(134, 13).Value = -4301.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 2923.2307
$ws.Cells.Item(132, 9).Value = 2875.5
$ws.Cells.Item(132, 10).Value = 2944.4443
$ws.Cells.Item(132, 11).Value = 25879.5
$ws.Cells.Item(132, 12).Value = 26499.9987
$ws.Cells.Item(132, 13).Value = -23349.5
$ws.Cells.Item(132, 14).Value = -31559.9987

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 127080.36
$ws.Cells.Item(70, 10).Value = 104586.57
$ws.Cells.Item(70, 12).Value = 104586.57
$ws.Cells.Item(70, 14).Value = -105126.57
$ws.Cells.Item(73, 8).Value = 127080.36
$ws.Cells.Item(73, 10).Value = 104586.57
$ws.Cells.Item(73, 12).Value = 104586.57
$ws.Cells.Item(73, 14).Value = -106458.57
$ws.Cells.Item(97, 8).Value = 500401.2
$ws.Cells.Item(97, 9).Value = 714608.2
$ws.Cells.Item(97, 11).Value = 714608.2
$ws.Cells.Item(97, 13).Value = -714112.2
$ws.Cells.Item(102, 8).Value = 1224.9584
$ws.Cells.Item(102, 9).Value = 1057.1904
$ws.Cells.Item(102, 10).Value = 2399.3333
$ws.Cells.Item(102, 11).Value = 1057.1904
$ws.Cells.Item(102, 12).Value = 2399.3333
$ws.Cells.Item(102, 13).Value = 564.8096
$ws.Cells.Item(102, 14).Value = -5643.3333
$ws.Cells.Item(113, 8).Value = 2651211.2
$ws.Cells.Item(113, 9).Value = 139800
$ws.Cells.Item(113, 10).Value = 6669469
$ws.Cells.Item(113, 11).Value = 139800
$ws.Cells.Item(113, 12).Value = 6669469
$ws.Cells.Item(113, 13).Value = -137630
$ws.Cells.Item(113, 14).Value = -6673809
$ws.Cells.Item(132, 8).Value = 5978.8887
$ws.Cells.Item(132, 9).Value = 5466.923
$ws.Cells.Item(132, 10).Value = 7310
$ws.Cells.Item(132, 11).Value = 16400.769
$ws.Cells.Item(132, 12).Value = 21930
$ws.Cells.Item(132, 13).Value = -13870.769
$ws.Cells.Item(132, 14).Value = -26990

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2077.5
$ws.Cells.Item(46, 10).Value = 3710
$ws.Cells.Item(46, 12).Value = 3710
$ws.Cells.Item(46, 14).Value = -4086
$ws.Cells.Item(68, 8).Value = 702333.3
$ws.Cells.Item(68, 9).Value = 702333.3
$ws.Cells.Item(68, 11).Value = 702333.3
$ws.Cells.Item(68, 13).Value = -701584.3
$ws.Cells.Item(71, 8).Value = 702333.3
$ws.Cells.Item(71, 9).Value = 702333.3
$ws.Cells.Item(71, 11).Value = 3511666.5
$ws.Cells.Item(71, 13).Value = -3507922.5
$ws.Cells.Item(82, 8).Value = 973.2857
$ws.Cells.Item(82, 9).Value = 966.1111
$ws.Cells.Item(82, 10).Value = 1016.3333
$ws.Cells.Item(82, 11).Value = 966.1111
$ws.Cells.Item(82, 12).Value = 1016.3333
$ws.Cells.Item(82, 13).Value = -605.1111
$ws.Cells.Item(82, 14).Value = -1738.3333
$ws.Cells.Item(85, 8).Value = 973.2857
$ws.Cells.Item(85, 9).Value = 966.1111
$ws.Cells.Item(85, 10).Value = 1016.3333
$ws.Cells.Item(85, 11).Value = 966.1111
$ws.Cells.Item(85, 12).Value = 1016.3333
$ws.Cells.Item(85, 13).Value = 281.8889
$ws.Cells.Item(85, 14).Value = -3512.3333
$ws.Cells.Item(97, 8).Value = 13200
$ws.Cells.Item(97, 10).Value = 13200
$ws.Cells.Item(97, 12).Value = 13200
$ws.Cells.Item(97, 14).Value = -15182
$ws.Cells.Item(108, 8).Value = 79999
$ws.Cells.Item(108, 10).Value = 79999
$ws.Cells.Item(108, 12).Value = 79999
$ws.Cells.Item(108, 14).Value = -87679
$ws.Cells.Item(122, 8).Value = 14320228
$ws.Cells.Item(122, 9).Value = 42926.637
$ws.Cells.Item(122, 11).Value = 128779.911
$ws.Cells.Item(122, 13).Value = -126329.911

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 49585.043
$ws.Cells.Item(126, 9).Value = 62696.945
$ws.Cells.Item(126, 11).Value = 188090.835
$ws.Cells.Item(126, 13).Value = -185620.835
$ws.Cells.Item(136, 8).Value = 1173.5
$ws.Cells.Item(136, 9).Value = 1157.7037
$ws.Cells.Item(136, 11).Value = 3473.1111
$ws.Cells.Item(136, 13).Value = -923.1111000000001
